$d = $word.ActiveDocument

$replacements = @(
    @{old="23×26="; new="46×49="},
    @{old="39×84="; new="44×18="},
    @{old="21×86="; new="57×57="},
    @{old="79×89="; new="20×18="},
    @{old="46×78="; new="22×24="},
    @{old="31×94="; new="96×26="},
    @{old="69×86="; new="46×37="},
    @{old="32×53="; new="28×95="},
    @{old="55×94="; new="29×72="},
    @{old="90×82="; new="61×54="},
    @{old="55×70="; new="89×69="},
    @{old="46×50="; new="48×48="},
    @{old="54×20="; new="56×31="},
    @{old="98×21="; new="19×80="},
    @{old="14×16="; new="43×58="},
    @{old="70×33="; new="42×37="},
    @{old="45×95="; new="85×21="},
    @{old="41×89="; new="30×42="},
    @{old="22×82="; new="37×77="},
    @{old="18×38="; new="18×93="},
    @{old="46×94="; new="82×61="},
    @{old="87×87="; new="34×55="},
    @{old="34×92="; new="82×60="},
    @{old="93×24="; new="55×47="},
    @{old="12×41="; new="73×61="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
